# Calibração pendulo.xlsx - add Sheet3 with the pendulum calibration curve,
# and tidy up the view/selection state on the existing sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Create Sheet3 as a copy of Sheet1 (keeps formulas/layout/quirky style
#     on C20 identical, then we overwrite the measured values). ---
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Sheet3"

# --- New calibration data for Sheet3 (A: angle, B/C: measured voltages) ---
$angles = @(0,10,20,30,40,50,60,70,80,90,100,110,120,130,140,150,160,170,180)
$v1 = @(2.91,2.5099999999999998,2.16,1.82,1.47,1.21,0.86,0.622,0.37,0.5,0.24,-0.54,-0.82,-1.05,-1.34,-1.58,-1.98,-2.27,-2.61)
$v2 = @(2.78,2.44,2.1800000000000002,1.7,1.39,1.1200000000000001,0.85,0.58199999999999996,0.32,0.7,0.25,-0.56999999999999995,-0.81,-1.06,-1.37,-1.57,-1.89,-2.21,-2.57)

for ($i = 0; $i -lt $angles.Count; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $angles[$i]
    $ws3.Cells.Item($row, 2).Value = $v1[$i]
    $ws3.Cells.Item($row, 3).Value = $v2[$i]
}

# C20 carried a (harmless, pre-existing) "quote prefix" cell style on Sheet1
# that a plain Value write clears; re-apply it from the untouched Sheet1
# cell so Sheet3's C20 keeps the same formatting as the source data.
$ws1.Range("C20").Copy()
$ws3.Range("C20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Recalculate so the CONCATENATE helper column (E) caches the new strings.
$excel.Calculate()

# --- View/selection tidy-up ---
# Sheet1: whole table selected, no single active cell highlighted.
$ws1.Range("A1:E20").Select()

# Sheet3: becomes the active tab, E2:E20 (the generated "{...}" column) selected.
$ws3.Range("E2:E20").Select()

Write-Host "done"
